$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 53, shifting existing rows 53:73 down to 54:74.
$ws.Rows("53:53").Insert()

# The newly inserted (blank) row 53 becomes a duplicate of the original row 52
# (this is what ends up as the new "last" row of the shifted block once the
# chain of weekly entries shifts down by one).
$ws.Range("A52:T52").Copy($ws.Range("A53:T53"))

# Row 52 itself gets this week's new entry: a new date and a new volume.
$ws.Range("D52").Value = 44460
$ws.Range("M52").Value = 200
